# The schedule sheets use a single-character "X" placeholder in the
# "grupa" (group) column of each day block as a stand-in value. This
# commit swaps every such literal "X" cell for a "-" across all sheets
# in the workbook.
$wb = $excel.ActiveWorkbook

$replaced = 0
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($cell.Value2 -eq "X") {
                $cell.Value2 = "-"
                $replaced++
            }
        }
    }
}

Write-Host "Replaced $replaced cell(s) containing 'X' with '-'."
